# Update view-count-style figures ("F" column) on several rows across the
# "展览" (sheet1) and "全部类型"/"本地生活" (sheet4/sheet3) sheets, matching
# the source data refresh from the "456a3b4" generated-output commit.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8188
$ws1.Range("F4").Value = 1918
$ws1.Range("F5").Value = 6522
$ws1.Range("F10").Value = 21
$ws1.Range("F16").Value = 8570
$ws1.Range("F33").Value = 2107
$ws1.Range("F37").Value = 2
$ws1.Range("F40").Value = 152
$ws1.Range("F41").Value = 6
$ws1.Range("F42").Value = 43

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 311

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8188
$ws4.Range("F7").Value = 311
$ws4.Range("F8").Value = 1918
$ws4.Range("F9").Value = 6522
$ws4.Range("F23").Value = 8570
$ws4.Range("F32").Value = 2107
